# feat: add 2022-Q1 data
#
# 1. Insert a new "2022-Q1" sheet (same layout/style as the other
#    quarterly sheets) right before the "总计" summary sheet.
# 2. Insert a new top row into "总计" summarising the new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, positioned right before "总计".
#
# Duplicate an existing quarterly sheet ("2021-Q1" - it already has 20
# pre-formatted data rows, plenty to cover the 10 rows we need) so the
# new tab inherits the workbook's sheetPr/pageMargins/header-row and
# index-column styling (bold, centered, thin borders) exactly, then
# wipe the copied values and write the real 2022-Q1 numbers.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q1")
$zongjiBefore = $wb.Worksheets.Item("总计")
$template.Copy($zongjiBefore)
$ws = $wb.Worksheets.Item("2021-Q1 (2)")
$ws.Name = "2022-Q1"

$ws.Cells.ClearContents()
$ws.Rows("12:21").Delete()

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'010723"
$ws.Range("C2").Value = "中欧价值成长混合A"
$ws.Range("D2").Value = "'28.80"
$ws.Range("E2").Value = "'91.96"
$ws.Range("F2").Value = "'3.77"
$ws.Range("G2").Value = "'1.0858"
$ws.Range("H2").Value = 9

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'166009"
$ws.Range("C3").Value = "中欧新动力混合(LOF) -A"
$ws.Range("D3").Value = "'24.29"
$ws.Range("E3").Value = "'90.17"
$ws.Range("F3").Value = "'3.81"
$ws.Range("G3").Value = "'0.9254"
$ws.Range("H3").Value = 10

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'001883"
$ws.Range("C4").Value = "中欧新动力混合(LOF) -E"
$ws.Range("D4").Value = "'24.29"
$ws.Range("E4").Value = "'90.17"
$ws.Range("F4").Value = "'3.81"
$ws.Range("G4").Value = "'0.9254"
$ws.Range("H4").Value = 10

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'010678"
$ws.Range("C5").Value = "中欧均衡成长混合A"
$ws.Range("D5").Value = "'18.02"
$ws.Range("E5").Value = "'90.22"
$ws.Range("F5").Value = "'3.69"
$ws.Range("G5").Value = "'0.6649"
$ws.Range("H5").Value = 9

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'009210"
$ws.Range("C6").Value = "中欧嘉和三年持有期混合A"
$ws.Range("D6").Value = "'18.38"
$ws.Range("E6").Value = "'85.26"
$ws.Range("F6").Value = "'3.61"
$ws.Range("G6").Value = "'0.6635"
$ws.Range("H6").Value = 8

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'005421"
$ws.Range("C7").Value = "中欧嘉泽灵活配置混合"
$ws.Range("D7").Value = "'8.87"
$ws.Range("E7").Value = "'86.45"
$ws.Range("F7").Value = "'3.12"
$ws.Range("G7").Value = "'0.2767"
$ws.Range("H7").Value = 10

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "'004236"
$ws.Range("C8").Value = "中欧新动力混合(LOF) -C"
$ws.Range("D8").Value = "'5.67"
$ws.Range("E8").Value = "'90.17"
$ws.Range("F8").Value = "'3.81"
$ws.Range("G8").Value = "'0.2160"
$ws.Range("H8").Value = 10

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "'009211"
$ws.Range("C9").Value = "中欧嘉和三年持有期混合C"
$ws.Range("D9").Value = "'2.42"
$ws.Range("E9").Value = "'85.26"
$ws.Range("F9").Value = "'3.61"
$ws.Range("G9").Value = "'0.0874"
$ws.Range("H9").Value = 8

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "'010724"
$ws.Range("C10").Value = "中欧价值成长混合C"
$ws.Range("D10").Value = "'1.88"
$ws.Range("E10").Value = "'91.96"
$ws.Range("F10").Value = "'3.77"
$ws.Range("G10").Value = "'0.0709"
$ws.Range("H10").Value = 9

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "'010679"
$ws.Range("C11").Value = "中欧均衡成长混合C"
$ws.Range("D11").Value = "'0.90"
$ws.Range("E11").Value = "'90.22"
$ws.Range("F11").Value = "'3.69"
$ws.Range("G11").Value = "'0.0332"
$ws.Range("H11").Value = 9

# ---------------------------------------------------------------------
# 2) "总计" sheet: insert the 2022-Q1 summary row at the top, pushing
#    the existing history rows down by one.
#
# NOTE: re-resolve the sheet by name here rather than reusing
# $zongjiBefore — worksheet references returned by Worksheets.Item
# track sheet *position*, and inserting the new "2022-Q1" tab above
# shifted "总计" from index 6 to index 7, so a handle captured before
# the Copy()/Add() now resolves to the wrong (new) tab.
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$zongji.Rows(2).Insert()
$zongji.Range("A2:D2").ClearFormats()
$zongji.Range("A3").Copy()
$zongji.Range("A2").PasteSpecial(-4122)

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q1"
$zongji.Range("C2").Value = 10
$zongji.Range("D2").Value = 4.95
